# Restructure the students sheet:
#  - rename "current credit hours" (F1) -> "Registered Credit Hours"
#  - rename "completed credit hours" (G1) -> "Earned Credit Hours"
#  - change the semester value (I2) from "SPRG" -> "Fall"
#  - add two new trailing columns: N "Attempted Credit hours" (66)
#    and O "Transcipt" (2)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "Earned Credit Hours"
$ws.Range("F1").Value = "Registered Credit Hours"

$ws.Range("N1").Value = "Attempted Credit hours"
$ws.Range("O1").Value = "Transcipt"

$ws.Range("I2").Value = "Fall"

$ws.Range("N2").Value = 66
$ws.Range("O2").Value = 2

$ws.Columns.Item(6).ColumnWidth = 24.45
$ws.Columns.Item(7).ColumnWidth = 24.17
$ws.Columns.Item(14).ColumnWidth = 22.65

$ws.Range("I2").Select()
